# process_04 inputData.xlsx update:
#  - add a new "Unitof measure" column (S) with value "lb" for the data row
#  - bump the "Valid From date" (H2) from 2021-02-15 to 2021-02-18
#  - move the sheet selection from O2 to H3 (scrolled so column B leads)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header + value
$ws.Range("S1").Value = "Unitof measure"
$ws.Range("S2").Value = "lb"

# Updated Valid From date
$ws.Range("H2").Value = (Get-Date -Year 2021 -Month 2 -Day 18).Date

# Update view state: scroll so column B is left-most visible, select H3
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("H3").Select()
